# Apply the changes described by the commit:
#  - Rename the "2014 Cohort" sheet to "2014_Cohort" (underscore instead of space)
#  - Clear the leftover numeric-format style (cellXfs index 3, applyNumberFormat)
#    that used to be stamped on every cell in column A (A1:A33)
#  - Remove the now-empty trailing cell A34 (it held no value, only that style)
#  - Move the active selection from D20 to B29

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the unused number-format styling that was applied down column A.
$ws.Range("A1:A33").ClearFormats()

# A34 had no value at all, just the stray style - clear it away entirely.
$ws.Range("A34").Clear()

# Rename the worksheet (space -> underscore).
$ws.Name = "2014_Cohort"

# Update the saved selection/active cell.
$ws.Range("B29").Select()
